# Updated cryptos list on Mon Feb 26 03:38:49 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to stay text even when the string looks like a number
    # (e.g. "384.41"), mirroring how Excel keeps a pre-existing text cell as
    # text instead of auto-converting the input to a numeric value.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "51.526.42"
$ws.Range("E2").Value = "  -0.06%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.099.49"
$ws.Range("E3").Value = "  +2.46%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.06%  "

# Row 5 - BNB
Set-TextValue "D5" "384.41"
$ws.Range("E5").Value = "  +1.40%  "

# Row 6 - Solana
$ws.Range("E6").Value = "  +0.30%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -0.67%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.06%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.45%  "

# Row 10 - Avalanche
Set-TextValue "D10" "36.91"
$ws.Range("E10").Value = "  +0.60%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.17%  "

# Row 12 - Dogecoin
$ws.Range("E12").Value = "  -0.12%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.585.31"
$ws.Range("E13").Value = "  +2.53%  "

# Row 14 - Chainlink
Set-TextValue "D14" "18.64"
$ws.Range("E14").Value = "  +0.81%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  +1.19%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.094.12"
$ws.Range("E16").Value = "  +2.62%  "

# Row 17 - Uniswap
Set-TextValue "D17" "11.15"
$ws.Range("E17").Value = "  +8.03%  "

# Row 18 - Polygon
$ws.Range("E18").Value = "  +1.03%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "51.534.44"
$ws.Range("E19").Value = "  -0.08%  "

# Row 20 - ImmutableX
Set-TextValue "D20" "3.33"
$ws.Range("E20").Value = "  +9.07%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "0.0₃0964"
$ws.Range("E21").Value = "  +0.31%  "

# Row 22 - InternetComputer(DFINITY)
$ws.Range("E22").Value = "  -0.65%  "

# Row 23 - Litecoin
Set-TextValue "D23" "69.96"
$ws.Range("E23").Value = "  -0.15%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "266.00"

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -0.59%  "

# Row 26 - Filecoin
$ws.Range("E26").Value = "  -0.91%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "27.04"
$ws.Range("E27").Value = "  +2.76%  "

# Row 28 - RenderToken
Set-TextValue "D28" "7.27"
$ws.Range("E28").Value = "  -2.47%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -1.56%  "

# Row 32 - Cosmos
$ws.Range("E32").Value = "  +0.18%  "

# Row 33 - InjectiveProtocol
Set-TextValue "D33" "35.35"
$ws.Range("E33").Value = "  +3.39%  "

# Row 34 - VeChain
Set-TextValue "D34" "0.0467"
$ws.Range("E34").Value = "  +3.17%  "

# Row 35 - Toncoin
Set-TextValue "D35" "2.07"
$ws.Range("E35").Value = "  +0.66%  "

# Row 36 - OKB
Set-TextValue "D36" "50.23"
$ws.Range("E36").Value = "  -0.67%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  -0.15%  "

# Row 38 - LidoDAOToken
Set-TextValue "D38" "3.37"
$ws.Range("E38").Value = "  +2.79%  "

# Row 39 - TheGraph
Set-TextValue "D39" "0.298"
$ws.Range("E39").Value = "  +5.25%  "

# Row 40 - ARBITRUM
Set-TextValue "D40" "1.89"
$ws.Range("E40").Value = "  +0.79%  "

# Row 41 - Monero
Set-TextValue "D41" "128.89"
$ws.Range("E41").Value = "  +1.22%  "

# Row 42 - was Stellar, now Celestia (swapped with row 43)
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D42" "16.57"
$ws.Range("E42").Value = "  -4.69%  "

# Row 43 - was Celestia, now Stellar (swapped with row 42)
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D43" "0.116"
$ws.Range("E43").Value = "  -0.21%  "

# Row 44 - Stacks
$ws.Range("E44").Value = "  -2.58%  "

# Row 45 - EnergySwap
Set-TextValue "D45" "22.44"
$ws.Range("E45").Value = "  +1.96%  "

# Row 46 - NEARProtocol
Set-TextValue "D46" "3.65"
$ws.Range("E46").Value = "  -2.01%  "

# Row 47 - ApeXProtocol
$ws.Range("E47").Value = "  +3.60%  "

# Row 48 - WEMIXToken
$ws.Range("E48").Value = "  +1.05%  "

# Row 49 - Maker
$ws.Range("D49").Value = "2.056.40"
$ws.Range("E49").Value = "  +1.37%  "

# Row 50 - BEAM
Set-TextValue "D50" "0.0329"
$ws.Range("E50").Value = "  +2.74%  "

# Row 51 - Mantle
Set-TextValue "D51" "0.893"
$ws.Range("E51").Value = "  +12.97%  "
